# TC04_Canine_Filter_StudyType-All.xlsx
# "FilesTab" row query (cell B4 on the "startup" sheet) is rewritten:
#  - the `File Type` output column is dropped
#  - the `Breed` output column is dropped
# (the rest of the Cypher query text is unchanged)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$newFileQuery = "MATCH (f:file)-->(parent)`n" + `
"WITH DISTINCT f, parent`n" + `
"MATCH (f)-[*]->(c:case)<--(demo:demographic)`n" + `
" MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)`n" + `
"WHERE s.clinical_study_type IN ['Clinical Trial','Transcriptomics','Genomics']`n" + `
"WITH DISTINCT f, parent, c, demo, diag, s`n" + `
"RETURN coalesce(f.file_name, '') AS ``File Name``, `n" + `
"        coalesce(labels(parent)[0], '') AS ``Association``,`n" + `
"        coalesce(f.file_description, '') AS ``Description``,`n" + `
"        coalesce(f.file_format, '') AS ``Format``,`n" + `
"        coalesce(f.file_size, '') AS ``Size``,`n" + `
"        coalesce(c.case_id, '') AS ``Case ID``, `n" + `
"        coalesce(diag.disease_term,'') AS Diagnosis , `n" + `
"        coalesce(s.clinical_study_designation,'') AS ``Study Code``"

$ws.Range("B4").Value = $newFileQuery

# reflect the author's final cursor position (scrolled down to row 4,
# cell B4 selected) on the worksheet
$ws.Activate()
$ws.Range("B4").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
